# Adaption of gap_types to 2 groups ("Arbeit" & "Privat") instead of
# 3 groups ("Arbeit", "Soziales Umfeld", "Haushalt & Selbstsorge").
#
# On sheet "QAGlist_Teil1" (the first/active sheet), column M ("Gap1_type")
# holds one of three group labels per question row. Every occurrence of
# "Haushalt & Selbstsorge" or "Soziales Umfeld" is consolidated into the
# single new label "Privat". Rows whose Gap1_type is already "Arbeit" stay
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QAGlist_Teil1")
$ws.Activate()

# Rows (in column M) that currently read "Haushalt & Selbstsorge" or
# "Soziales Umfeld" and must become "Privat".
$rowsToUpdate = @(3, 4, 5, 6, 7, 8, 9, 10, 19, 22, 24, 25, 26, 28, 29, 30)

foreach ($r in $rowsToUpdate) {
    $ws.Cells.Item($r, 13).Value = "Privat"
}

# Reflect the view state recorded in the saved workbook: the window is
# scrolled so column C is left-most, and M30 is the active/selected cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("M30").Select()
